# Fix the "Recorded By" (column G) values: the author/System order was
# swapped from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for a specific subset of the attendance rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "dnasr281@gmail.com, System"
$ws.Range("G3").Value = "dnasr281@gmail.com, System"
$ws.Range("G4").Value = "dnasr281@gmail.com, System"
$ws.Range("G5").Value = "dnasr281@gmail.com, System"
$ws.Range("G6").Value = "dnasr281@gmail.com, System"
$ws.Range("G7").Value = "dnasr281@gmail.com, System"
$ws.Range("G16").Value = "dnasr281@gmail.com, System"
$ws.Range("G17").Value = "dnasr281@gmail.com, System"
$ws.Range("G22").Value = "dnasr281@gmail.com, System"
$ws.Range("G23").Value = "dnasr281@gmail.com, System"
$ws.Range("G37").Value = "dnasr281@gmail.com, System"
$ws.Range("G38").Value = "dnasr281@gmail.com, System"
$ws.Range("G43").Value = "dnasr281@gmail.com, System"
$ws.Range("G44").Value = "dnasr281@gmail.com, System"
$ws.Range("G58").Value = "dnasr281@gmail.com, System"
$ws.Range("G59").Value = "dnasr281@gmail.com, System"
$ws.Range("G64").Value = "dnasr281@gmail.com, System"
$ws.Range("G65").Value = "dnasr281@gmail.com, System"
$ws.Range("G79").Value = "dnasr281@gmail.com, System"
$ws.Range("G80").Value = "dnasr281@gmail.com, System"
$ws.Range("G85").Value = "dnasr281@gmail.com, System"
$ws.Range("G86").Value = "dnasr281@gmail.com, System"
$ws.Range("G87").Value = "dnasr281@gmail.com, System"
$ws.Range("G88").Value = "dnasr281@gmail.com, System"
$ws.Range("G89").Value = "dnasr281@gmail.com, System"
$ws.Range("G90").Value = "dnasr281@gmail.com, System"
$ws.Range("G99").Value = "dnasr281@gmail.com, System"
$ws.Range("G100").Value = "dnasr281@gmail.com, System"
$ws.Range("G105").Value = "dnasr281@gmail.com, System"
$ws.Range("G106").Value = "dnasr281@gmail.com, System"
$ws.Range("G107").Value = "dnasr281@gmail.com, System"
$ws.Range("G108").Value = "dnasr281@gmail.com, System"
$ws.Range("G109").Value = "dnasr281@gmail.com, System"
$ws.Range("G110").Value = "dnasr281@gmail.com, System"
$ws.Range("G119").Value = "dnasr281@gmail.com, System"
$ws.Range("G120").Value = "dnasr281@gmail.com, System"
$ws.Range("G125").Value = "dnasr281@gmail.com, System"
$ws.Range("G126").Value = "dnasr281@gmail.com, System"
$ws.Range("G127").Value = "dnasr281@gmail.com, System"
$ws.Range("G128").Value = "dnasr281@gmail.com, System"
$ws.Range("G129").Value = "dnasr281@gmail.com, System"
$ws.Range("G130").Value = "dnasr281@gmail.com, System"
$ws.Range("G139").Value = "dnasr281@gmail.com, System"
$ws.Range("G140").Value = "dnasr281@gmail.com, System"
$ws.Range("G145").Value = "dnasr281@gmail.com, System"
$ws.Range("G146").Value = "dnasr281@gmail.com, System"
$ws.Range("G147").Value = "dnasr281@gmail.com, System"
$ws.Range("G148").Value = "dnasr281@gmail.com, System"
$ws.Range("G149").Value = "dnasr281@gmail.com, System"
$ws.Range("G150").Value = "dnasr281@gmail.com, System"
$ws.Range("G159").Value = "dnasr281@gmail.com, System"
$ws.Range("G160").Value = "dnasr281@gmail.com, System"
$ws.Range("G165").Value = "dnasr281@gmail.com, System"
$ws.Range("G166").Value = "dnasr281@gmail.com, System"
$ws.Range("G167").Value = "dnasr281@gmail.com, System"
$ws.Range("G168").Value = "dnasr281@gmail.com, System"
$ws.Range("G169").Value = "dnasr281@gmail.com, System"
$ws.Range("G170").Value = "dnasr281@gmail.com, System"
$ws.Range("G179").Value = "dnasr281@gmail.com, System"
$ws.Range("G180").Value = "dnasr281@gmail.com, System"
$ws.Range("G185").Value = "dnasr281@gmail.com, System"
$ws.Range("G186").Value = "dnasr281@gmail.com, System"
$ws.Range("G200").Value = "dnasr281@gmail.com, System"
$ws.Range("G201").Value = "dnasr281@gmail.com, System"
$ws.Range("G206").Value = "dnasr281@gmail.com, System"
$ws.Range("G207").Value = "dnasr281@gmail.com, System"
$ws.Range("G221").Value = "dnasr281@gmail.com, System"
$ws.Range("G222").Value = "dnasr281@gmail.com, System"
$ws.Range("G227").Value = "dnasr281@gmail.com, System"
$ws.Range("G228").Value = "dnasr281@gmail.com, System"
$ws.Range("G242").Value = "dnasr281@gmail.com, System"
$ws.Range("G243").Value = "dnasr281@gmail.com, System"
